$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TCLRY")

$ws.Range("D8").Value = 1984800
$ws.Range("E8").Value = 2339300
$ws.Range("F8").Value = 2353900
$ws.Range("G8").Value = 2477400
$ws.Range("H8").Value = 1273500
$ws.Range("I8").Value = 2715200
$ws.Range("J8").Value = 1416000
$ws.Range("D9").Value = 1762600
$ws.Range("E9").Value = 2022900
$ws.Range("F9").Value = 2055500
$ws.Range("G9").Value = 2205800
$ws.Range("I9").Value = 2209200
$ws.Range("D10").Value = 222200
$ws.Range("E10").Value = 316400
$ws.Range("F10").Value = 298400
$ws.Range("G10").Value = 271500
$ws.Range("I10").Value = 506000
$ws.Range("D12").Value = 86400
$ws.Range("E12").Value = 95400
$ws.Range("F12").Value = 97600
$ws.Range("G12").Value = 94200
$ws.Range("I12").Value = 104300
$ws.Range("D14").Value = 46000
$ws.Range("E14").Value = 29200
$ws.Range("F14").Value = 29200
$ws.Range("G14").Value = 5600
$ws.Range("I14").Value = 52700
$ws.Range("D17").Value = 2103700
$ws.Range("E17").Value = 2309100
$ws.Range("F17").Value = 2417900
$ws.Range("G17").Value = 2498700
$ws.Range("I17").Value = 2608600
$ws.Range("D18").Value = -118900
$ws.Range("E18").Value = 30300
$ws.Range("F18").Value = -64000
$ws.Range("G18").Value = -21300
$ws.Range("I18").Value = 106600
$ws.Range("E20").Value = -15700
$ws.Range("F20").Value = -41500
$ws.Range("G20").Value = -43800
$ws.Range("I20").Value = -31400
$ws.Range("D21").Value = 3400
$ws.Range("E21").Value = 154800
$ws.Range("F21").Value = 23600
$ws.Range("G21").Value = 68400
$ws.Range("I21").Value = 191900
$ws.Range("D22").Value = 23600
$ws.Range("E22").Value = 23600
$ws.Range("F22").Value = 28000
$ws.Range("G22").Value = 44900
$ws.Range("I22").Value = 50500
$ws.Range("D23").Value = -142500
$ws.Range("E23").Value = -9000
$ws.Range("F23").Value = -133500
$ws.Range("G23").Value = -110000
$ws.Range("I23").Value = 24700
$ws.Range("D24").Value = 12300
$ws.Range("E24").Value = 113300
$ws.Range("F24").Value = 6700
$ws.Range("I24").Value = 33700
$ws.Range("D26").Value = -154800
$ws.Range("E26").Value = -122300
$ws.Range("F26").Value = -140200
$ws.Range("G26").Value = -110000
$ws.Range("I26").Value = -9000
$ws.Range("D27").Value = -154800
$ws.Range("E27").Value = -122300
$ws.Range("F27").Value = -140200
$ws.Range("G27").Value = -110000
$ws.Range("I27").Value = -9000
$ws.Range("D29").Value = -15700
$ws.Range("E29").Value = 47100
$ws.Range("F29").Value = 21300
$ws.Range("G29").Value = 139100
$ws.Range("I29").Value = -49400
$ws.Range("E32").Value = 15700
$ws.Range("F32").Value = 41500
$ws.Range("G32").Value = 43800
$ws.Range("I32").Value = 31400
$ws.Range("D33").Value = -170500
$ws.Range("E33").Value = -75200
$ws.Range("F33").Value = -118900
$ws.Range("G33").Value = 29200
$ws.Range("I33").Value = -58300
$ws.Range("D35").Value = -170500
$ws.Range("E35").Value = -75200
$ws.Range("F35").Value = -118900
$ws.Range("G35").Value = 29200
$ws.Range("I35").Value = -58300
$ws.Range("D41").Value = 221000
$ws.Range("E41").Value = 357900
$ws.Range("F41").Value = 205300
$ws.Range("G41").Value = 416300
$ws.Range("I41").Value = 486900
$ws.Range("D42").Value = 10100
$ws.Range("E42").Value = 11200
$ws.Range("F42").Value = 13500
$ws.Range("G42").Value = 19100
$ws.Range("I42").Value = 25800
$ws.Range("D43").Value = 759600
$ws.Range("E43").Value = 809000
$ws.Range("F43").Value = 892000
$ws.Range("G43").Value = 963800
$ws.Range("I43").Value = 892000
$ws.Range("D44").Value = 271500
$ws.Range("E44").Value = 267000
$ws.Range("F44").Value = 267000
$ws.Range("G44").Value = 262500
$ws.Range("I44").Value = 245700
$ws.Range("D45").Value = 233400
$ws.Range("E45").Value = 295100
$ws.Range("F45").Value = 268200
$ws.Range("G45").Value = 318600
$ws.Range("I45").Value = 482500
$ws.Range("D46").Value = 1495600
$ws.Range("E46").Value = 1740200
$ws.Range("F46").Value = 1646000
$ws.Range("G46").Value = 1980300
$ws.Range("I46").Value = 2132900
$ws.Range("D47").Value = 38100
$ws.Range("E47").Value = 42600
$ws.Range("F47").Value = 46000
$ws.Range("G47").Value = 68400
$ws.Range("I47").Value = 80800
$ws.Range("D48").Value = 250200
$ws.Range("E48").Value = 272600
$ws.Range("F48").Value = 288400
$ws.Range("G48").Value = 320900
$ws.Range("I48").Value = 325400
$ws.Range("D49").Value = 1773900
$ws.Range("E49").Value = 1758200
$ws.Range("F49").Value = 1868100
$ws.Range("G49").Value = 2008400
$ws.Range("I49").Value = 1938800
$ws.Range("D52").Value = 360200
$ws.Range("E52").Value = 351200
$ws.Range("F52").Value = 499300
$ws.Range("G52").Value = 537400
$ws.Range("I52").Value = 578900
$ws.Range("D54").Value = 3918000
$ws.Range("E54").Value = 4164800
$ws.Range("F54").Value = 4347700
$ws.Range("G54").Value = 4915400
$ws.Range("I54").Value = 5056800
$ws.Range("D57").Value = 953700
$ws.Range("E57").Value = 1062500
$ws.Range("F57").Value = 999700
$ws.Range("G57").Value = 1113000
$ws.Range("I57").Value = 858300
$ws.Range("D58").Value = 24700
$ws.Range("E58").Value = 22400
$ws.Range("F58").Value = 13500
$ws.Range("G58").Value = 58300
$ws.Range("I58").Value = 96500
$ws.Range("D59").Value = 823500
$ws.Range("E59").Value = 787600
$ws.Range("F59").Value = 793200
$ws.Range("G59").Value = 958200
$ws.Range("I59").Value = 1147800
$ws.Range("D60").Value = 1801900
$ws.Range("E60").Value = 1872600
$ws.Range("F60").Value = 1806400
$ws.Range("G60").Value = 2129500
$ws.Range("I60").Value = 2102600
$ws.Range("D61").Value = 1217400
$ws.Range("E61").Value = 1208400
$ws.Range("F61").Value = 1211700
$ws.Range("G61").Value = 1119700
$ws.Range("I61").Value = 1320600
$ws.Range("D62").Value = 691100
$ws.Range("E62").Value = 706900
$ws.Range("F62").Value = 786500
$ws.Range("G62").Value = 876300
$ws.Range("I62").Value = 941300
$ws.Range("D66").Value = 3713800
$ws.Range("E66").Value = 3791200
$ws.Range("F66").Value = 3808000
$ws.Range("G66").Value = 4128900
$ws.Range("I66").Value = 4366800
$ws.Range("D72").Value = -1571900
$ws.Range("E72").Value = -1401400
$ws.Range("F72").Value = -1293700
$ws.Range("G72").Value = -1164600
$ws.Range("I72").Value = -1240900
$ws.Range("D76").Value = 204200
$ws.Range("E76").Value = 373600
$ws.Range("F76").Value = 539700
$ws.Range("G76").Value = 786500
$ws.Range("I76").Value = 690000
$ws.Range("D81").Value = -170500
$ws.Range("E81").Value = -75200
$ws.Range("F81").Value = -118900
$ws.Range("G81").Value = 29200
$ws.Range("I81").Value = -58300
$ws.Range("D83").Value = 122300
$ws.Range("E83").Value = 135800
$ws.Range("F83").Value = 133500
$ws.Range("G83").Value = 142500
$ws.Range("I83").Value = 116700
$ws.Range("D89").Value = -50500
$ws.Range("E89").Value = 274900
$ws.Range("F89").Value = -88600
$ws.Range("G89").Value = 264800
$ws.Range("I89").Value = 193000
$ws.Range("D91").Value = -33700
$ws.Range("E91").Value = -30300
$ws.Range("F91").Value = -28000
$ws.Range("G91").Value = -37000
$ws.Range("I91").Value = -39300
$ws.Range("D94").Value = -72900
$ws.Range("E94").Value = -88600
$ws.Range("F94").Value = -80800
$ws.Range("G94").Value = -46000
$ws.Range("I94").Value = -84100
$ws.Range("F96").Value = -28000
$ws.Range("I96").Value = -28000
$ws.Range("D100").Value = 2200
$ws.Range("E100").Value = -21300
$ws.Range("F100").Value = -10100
$ws.Range("G100").Value = -305200
$ws.Range("I100").Value = -58300
$ws.Range("D101").Value = -15700
$ws.Range("E101").Value = -19100
$ws.Range("F101").Value = -24700
$ws.Range("G101").Value = 24700
$ws.Range("I101").Value = 4500
$ws.Range("D102").Value = -136900
$ws.Range("E102").Value = 152600
$ws.Range("F102").Value = -210900
$ws.Range("G102").Value = -70700
$ws.Range("I102").Value = 55000
